{"js": "// Set Arial as the default font for the document (Normal style) and for\n// all heading styles (Heading 1 - Heading 9), mirroring the commit's\n// \"_set_default_font()\" behaviour:\n//   - Normal style: Arial font, 10pt (sz=20 half-points)\n//   - Heading 1-9 styles: add Arial as the ascii/hAnsi font, keeping the\n//     existing theme font references (majorHAnsi/majorEastAsia/majorBidi)\n//     and sizes/colors untouched.\n\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\nconst targetNames = [\n  \"Normal\",\n  \"Heading 1\",\n  \"Heading 2\",\n  \"Heading 3\",\n  \"Heading 4\",\n  \"Heading 5\",\n  \"Heading 6\",\n  \"Heading 7\",\n  \"Heading 8\",\n  \"Heading 9\",\n];\n\nfor (const style of styles.items) {\n  if (!targetNames.includes(style.nameLocal)) {\n    continue;\n  }\n\n  // Arial ascii/hAnsi font for every targeted style.\n  style.font.name = \"Arial\";\n\n  // Normal is the document's base style; also pin its size to 10pt so the\n  // whole document defaults to Arial 10pt.\n  if (style.nameLocal === \"Normal\") {\n    style.font.size = 10;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Set Arial as the default font for the document (Normal style) and for\n# all heading styles (Heading 1 - Heading 9), mirroring the commit's\n# \"_set_default_font()\" behaviour:\n#   - Normal style: Arial font, 10pt (sz=20 half-points)\n#   - Heading 1-9 styles: add Arial as the ascii/hAnsi font, keeping the\n#     existing theme font references (majorHAnsi/majorEastAsia/majorBidi)\n#     and sizes/colors untouched.\n\n$d = $word.ActiveDocument\n\n$targetNames = @(\n    \"Normal\",\n    \"Heading 1\",\n    \"Heading 2\",\n    \"Heading 3\",\n    \"Heading 4\",\n    \"Heading 5\",\n    \"Heading 6\",\n    \"Heading 7\",\n    \"Heading 8\",\n    \"Heading 9\"\n)\n\nforeach ($s in $d.Styles) {\n    if ($targetNames -contains $s.NameLocal) {\n        $s.Font.Name = \"Arial\"\n\n        if ($s.NameLocal -eq \"Normal\") {\n            $s.Font.Size = 10\n        }\n    }\n}\n"}
